# Remove Plasma Cell markers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 7, 8, 9 hold the "Plasma cells" entries (Plasma cells/MZB1, Plasma cells/IGHE, Plasma cells/IGKC).
# Delete that 3-row block; rows below shift up automatically.
$ws.Range("A7:B9").EntireRow.Delete() | Out-Null

# Reflect the selection left behind in the saved file (selecting the rows that used to
# contain the Plasma cells block, now occupied by the shifted-up data).
$ws.Range("A7:XFD9").Select() | Out-Null
